$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the FFR column (column C) and the FFR Lag row (row 3)
$ws.Range("C1:C4").EntireColumn.Delete() | Out-Null
$ws.Range("A3:C3").EntireRow.Delete() | Out-Null

# Update the remaining values with the new data (force text so numeric-looking
# values like "0.27" are stored as strings, matching the other text entries)
$ws.Range("B2").Value = "'0.27"
$ws.Range("B3").Value = "'-0.11*"
$ws.Range("C2").Value = "'-8.77*"
$ws.Range("C3").Value = "'2.45***"
$ws.Range("B2:C3").Style = "Normal"
